{"js": "// Update published course materials: Module 09 title + learning objectives\n// rewritten from \"Mendelian Genetics\" topics to \"Tissues and the Animal Body\" topics.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst replacements = [\n  {\n    find: \"Module 09: Tissues \u2014 Keys to Success\",\n    replace: \"Module 09: Tissues and the Animal Body \u2014 Keys to Success\",\n  },\n  {\n    find: \"1. Mendelian Genetics Foundations\",\n    replace: \"1. Homeostasis and Osmoregulation\",\n  },\n  {\n    find: \"Define gene, allele, genotype, and phenotype  Distinguish between dominant and recessive alleles  Explain homozygous and heterozygous conditions  Describe Mendel's experiments and his two laws   2. Law of Segregation\",\n    replace: \"Define homeostasis and describe its role in maintaining a stable internal environment  Explain negative and positive feedback loops using biological examples  Define osmoregulation and explain how the body balances water and salt   2. Digestive System\",\n  },\n  {\n    find: \"Explain the law of segregation  Use Punnett squares to predict offspring genotypes and phenotypes  Calculate genotypic and phenotypic ratios for monohybrid crosses   3. Law of Independent Assortment\",\n    replace: \"Describe the primary function of the digestive system  Identify the main organs of the digestive tract and their specific roles (mouth, stomach, small intestine, large intestine)  Explain the difference between mechanical and chemical digestion  Describe the role of accessory organs (liver, pancreas, gallbladder)   3. Circulatory and Respiratory Systems\",\n  },\n  {\n    find: \"Explain the law of independent assortment  Perform dihybrid crosses using Punnett squares  Understand when independent assortment applies and its limitations   4. Extensions to Mendelian Genetics\",\n    replace: \"Explain the function of the circulatory system in transporting nutrients and waste  Describe the pathway of blood through the human heart and lungs  Differentiate between arteries, veins, and capillaries  Explain how the respiratory system facilitates gas exchange (oxygen and carbon dioxide)   4. Endocrine System\",\n  },\n  {\n    find: \"Describe incomplete dominance and codominance  Explain multiple alleles using ABO blood types as an example  Understand polygenic inheritance and continuous variation  Explain pleiotropy (one gene affecting multiple traits)   5. Sex Linkage and Chromosomal Inheritance\",\n    replace: \"Define hormones and describe how the endocrine system uses them to communicate  Identify major endocrine glands (pituitary, thyroid, adrenal, pancreas) and their primary functions  Explain how the endocrine system works with the nervous system to maintain homeostasis   5. Musculoskeletal System\",\n  },\n  {\n    find: \"Describe sex determination in humans  Explain sex-linked inheritance patterns  Predict outcomes of crosses involving X-linked traits  Understand why sex-linked disorders are more common in males   6. Pedigree Analysis\",\n    replace: \"Describe the primary functions of the skeletal system (support, protection, movement, blood cell production)  Differentiate between the axial and appendicular skeleton  Differentiate between skeletal, smooth, and cardiac muscle tissue  Explain how muscles and bones work together to create movement   6. Nervous System\",\n  },\n  {\n    find: \"Interpret pedigree charts  Determine modes of inheritance from pedigrees  Identify carriers and affected individuals    Study Tips\",\n    replace: \"Describe the primary function of the nervous system in processing information  Differentiate between the central nervous system (CNS) and peripheral nervous system (PNS)  Identify the basic structure of a neuron and describe how signals are transmitted    Study Tips\",\n  },\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const currentText = para.text;\n  for (const { find, replace } of replacements) {\n    if (currentText === find) {\n      para.insertText(replace, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update published course materials: Module 09 title + learning objectives\n# rewritten from \"Mendelian Genetics\" topics to \"Tissues and the Animal Body\" topics.\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{\n        Find = \"Module 09: Tissues \u2014 Keys to Success\"\n        Replace = \"Module 09: Tissues and the Animal Body \u2014 Keys to Success\"\n    },\n    @{\n        Find = \"1. Mendelian Genetics Foundations\"\n        Replace = \"1. Homeostasis and Osmoregulation\"\n    },\n    @{\n        Find = \"Define gene, allele, genotype, and phenotype  Distinguish between dominant and recessive alleles  Explain homozygous and heterozygous conditions  Describe Mendel's experiments and his two laws   2. Law of Segregation\"\n        Replace = \"Define homeostasis and describe its role in maintaining a stable internal environment  Explain negative and positive feedback loops using biological examples  Define osmoregulation and explain how the body balances water and salt   2. Digestive System\"\n    },\n    @{\n        Find = \"Explain the law of segregation  Use Punnett squares to predict offspring genotypes and phenotypes  Calculate genotypic and phenotypic ratios for monohybrid crosses   3. Law of Independent Assortment\"\n        Replace = \"Describe the primary function of the digestive system  Identify the main organs of the digestive tract and their specific roles (mouth, stomach, small intestine, large intestine)  Explain the difference between mechanical and chemical digestion  Describe the role of accessory organs (liver, pancreas, gallbladder)   3. Circulatory and Respiratory Systems\"\n    },\n    @{\n        Find = \"Explain the law of independent assortment  Perform dihybrid crosses using Punnett squares  Understand when independent assortment applies and its limitations   4. Extensions to Mendelian Genetics\"\n        Replace = \"Explain the function of the circulatory system in transporting nutrients and waste  Describe the pathway of blood through the human heart and lungs  Differentiate between arteries, veins, and capillaries  Explain how the respiratory system facilitates gas exchange (oxygen and carbon dioxide)   4. Endocrine System\"\n    },\n    @{\n        Find = \"Describe incomplete dominance and codominance  Explain multiple alleles using ABO blood types as an example  Understand polygenic inheritance and continuous variation  Explain pleiotropy (one gene affecting multiple traits)   5. Sex Linkage and Chromosomal Inheritance\"\n        Replace = \"Define hormones and describe how the endocrine system uses them to communicate  Identify major endocrine glands (pituitary, thyroid, adrenal, pancreas) and their primary functions  Explain how the endocrine system works with the nervous system to maintain homeostasis   5. Musculoskeletal System\"\n    },\n    @{\n        Find = \"Describe sex determination in humans  Explain sex-linked inheritance patterns  Predict outcomes of crosses involving X-linked traits  Understand why sex-linked disorders are more common in males   6. Pedigree Analysis\"\n        Replace = \"Describe the primary functions of the skeletal system (support, protection, movement, blood cell production)  Differentiate between the axial and appendicular skeleton  Differentiate between skeletal, smooth, and cardiac muscle tissue  Explain how muscles and bones work together to create movement   6. Nervous System\"\n    },\n    @{\n        Find = \"Interpret pedigree charts  Determine modes of inheritance from pedigrees  Identify carriers and affected individuals    Study Tips\"\n        Replace = \"Describe the primary function of the nervous system in processing information  Differentiate between the central nervous system (CNS) and peripheral nervous system (PNS)  Identify the basic structure of a neuron and describe how signals are transmitted    Study Tips\"\n    }\n)\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $item.Find\n    $find.Replacement.Text = $item.Replace\n    $find.Execute($item.Find, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $item.Replace, $wdReplaceAll)\n}\n"}
